# Refresh cryptos list: updates Price (D) and Volume(1h) (E) columns,
# and corrects two coin rows whose rank order swapped (B/C/D/E).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "69.830.39"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "3.704.99"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'678.33"
$ws.Range("E5").Value = "  -1.24%  "
$ws.Range("D6").Value = "'162.68"
$ws.Range("E6").Value = "  +1.70%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +0.81%  "
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("D10").Value = "'7.13"
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("E11").Value = "  +2.17%  "
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("D13").Value = "'32.96"
$ws.Range("E13").Value = "  +1.46%  "
$ws.Range("D14").Value = "3.719.20"
$ws.Range("D15").Value = "69.789.82"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("E16").Value = "  +1.70%  "
$ws.Range("D17").Value = "'16.14"
$ws.Range("E17").Value = "  +1.93%  "
$ws.Range("D18").Value = "'6.51"
$ws.Range("E18").Value = "  +1.55%  "
$ws.Range("D19").Value = "'473.05"
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("D20").Value = "'9.84"
$ws.Range("E20").Value = "  -1.34%  "
$ws.Range("D21").Value = "'0.654"
$ws.Range("E21").Value = "  +0.58%  "
$ws.Range("D22").Value = "'80.58"
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("D23").Value = "3.853.96"
$ws.Range("E24").Value = "  +3.20%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "'11.04"
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").Value = "'9.16"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("E29").Value = "  +0.63%  "
$ws.Range("D30").Value = "'2.03"
$ws.Range("E30").Value = "  +1.49%  "
$ws.Range("D31").Value = "'6.64"
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'26.94"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").Value = "'0.165"
$ws.Range("E34").Value = "  +3.69%  "
$ws.Range("D35").Value = "3.694.98"
$ws.Range("E35").Value = "  +0.85%  "
$ws.Range("D36").Value = "'8.60"
$ws.Range("E36").Value = "  +4.63%  "
$ws.Range("E37").Value = "  +1.07%  "
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'2.23"
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "'0.0909"
$ws.Range("E41").Value = "  +1.10%  "
$ws.Range("D42").Value = "'168.88"
$ws.Range("E42").Value = "  +2.00%  "
$ws.Range("D43").Value = "'0.946"
$ws.Range("E43").Value = "  +0.26%  "
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("E45").Value = "  +2.36%  "
$ws.Range("E46").Value = "  -0.87%  "
$ws.Range("D47").Value = "'28.15"
$ws.Range("E47").Value = "  +2.28%  "
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("E50").Value = "  +2.10%  "
$ws.Range("E51").Value = "  +2.65%  "
